$wb = $excel.ActiveWorkbook

# The workbook has two sheets with the same "想去人数" (F column) data that
# need to be refreshed with newly scraped counts: "展览" and "全部类型".
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Map of row number -> new value for column F, keyed by the event name
    # in column C so we update the correct row regardless of row offset
    # differences between sheets.
    $updates = @{
        "南宁·第二届北极光动漫展"              = 3440
        "南宁·原神x星铁x绝区零同人ONLY3.0"        = 70
        "南宁·2024良牙动漫秋季盛典（秋典）"        = 1729
        "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini" = 91
        "南宁·万圣漫控嘉年华10"                  = 337
    }

    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value()
        if ($updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
